$wb = $excel.ActiveWorkbook

# 1. Copy the "Outofstock" sheet and place it right after "DataSet";
#    this becomes the new "Search" results-validation sheet (sheetId 8,
#    same shared header row / hyperlinked cells as Outofstock rows 1-3).
$src = $wb.Worksheets.Item("Outofstock")
$after = $wb.Worksheets.Item("DataSet")
$src.Copy([System.Reflection.Missing]::Value, $after)
$ws = $wb.Worksheets.Item(2)
$ws.Name = "Search"

# 2. Extra header cells (Y1:AC1) get the same yellow header style as the
#    rest of row 1.
$ws.Range("X1").Copy()
$ws.Range("Y1:AC1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("Y1").Value = "Products"
$ws.Range("Z1").Value = "Quantity"
$ws.Range("AA1").Value = "Invalid_searchdata"
$ws.Range("AB1").Value = "Sortby_Dropdown"
$ws.Range("AC1").Value = "Price_Symbol"

# 3. New validation rows 4-7 - written in the same order the values were
#    first introduced so shared-string indices line up.
$ws.Range("A4").Value = "Invalid_Search"

$ws.Range("Y4").Value = "q@!e#d`$D"
$ws.Hyperlinks.Add($ws.Range("Y4"), "mailto:q@!e#d`$D")
$ws.Range("B2").Copy()
$ws.Range("Y4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A5").Value = "SortBy"
$ws.Range("AB5").Value = "Highest Price"
$ws.Range("AC5").Value = "$"

$ws.Range("A6").Value = "Valid_Search"
$ws.Range("K6").Style = "Hyperlink"
$ws.Range("L6").Style = "Hyperlink"
$ws.Range("Y6").Value = "Hot Toddy Heat Protectant Mist"

$ws.Range("M2").Copy()
$ws.Range("Z6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Z6").Value = "'1"

$ws.Range("A7").Value = "Product"

# 4. Column widths for the new columns.
$ws.Range("Y1").ColumnWidth = 19.92
$ws.Range("Z1:AA1").ColumnWidth = 15.6
$ws.Range("AB1").ColumnWidth = 10.6

# 5. View state - "Search" becomes the active/selected tab; "DataSet"
#    keeps a plain (non-selected) view scrolled/selected elsewhere.
$dataSet = $wb.Worksheets.Item("DataSet")
$dataSet.Activate()
$dataSet.Range("K42").Select()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 22

$ws.Activate()
$ws.Range("J13").Select()
